$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained one new weekly observation. Insert a new row at 326,
# pushing the previous rows 326-357 down to 327-358 (matching the diff,
# which is a uniform "row 326 data downward shift" plus one brand-new row).
$ws.Rows(326).Insert()

# Populate the new row 326 with the new "Ajo" (garlic) price observation.
$ws.Cells.Item(326, 1).Value = 9
$ws.Cells.Item(326, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(326, 3).Value = "Metropolitana"
$ws.Cells.Item(326, 4).Value = 45194
$ws.Cells.Item(326, 5).Value = 13
$ws.Cells.Item(326, 6).Value = 100112003
$ws.Cells.Item(326, 7).Value = "Ajo"
$ws.Cells.Item(326, 8).Value = "Chino"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 520
$ws.Cells.Item(326, 11).Value = 19000
$ws.Cells.Item(326, 12).Value = 20000
$ws.Cells.Item(326, 13).Value = 19500
$ws.Cells.Item(326, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(326, 15).Value = "China"
$ws.Cells.Item(326, 16).Value = 1950
$ws.Cells.Item(326, 17).Value = 10
$ws.Cells.Item(326, 18).Value = "Hortaliza"
